$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.747.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.287.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.57%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.285.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.53%  "

$ws.Range("E10").Value = "  -5.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.08%  "

$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("E13").Value = "  -4.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.697.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.769.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.298.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.64%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.18%  "

$ws.Range("E25").Value = "  -4.48%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  -6.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("E30").Value = "  -6.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0717"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.48%  "

$ws.Range("E32").Value = "  -6.09%  "

$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.377"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.07%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -7.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.52%  "

$ws.Range("E44").Value = "  -4.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.51%  "

$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.554"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.12%  "

$ws.Range("E49").Value = "  -4.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("E51").Value = "  +82.85%  "
